$wb = $excel.ActiveWorkbook
$wsIn = $wb.Worksheets.Item("ProductLoan_Input")
$wsOut = $wb.Worksheets.Item("ProductLoan_Output")

# --- ProductLoan_Input sheet edits ---

# productname changes from "...RBI-EPP-..." to "486-RBI-EPP-..."
$wsIn.Range("B1").Value = "486-RBI-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"

# description text loses its leading newline
$wsIn.Range("B2").Value = "RBI (India) - Equal principle Installment - Declining Balance - Daily-Enable Check box Recalculate Interest-Interest recalculation compounding on NONE - Advance Reduce number of installments- Pre Calculate Till preclosure date- frequency for recalculate Outstanding Principal Daily - Checkbox Enable Multiple Disbursals yes - Maximum Tranche count 1"

# shortname changes from text "kar3" to numeric 486
$wsIn.Range("B3").Value = 486

# nominalinterestratedefault changes from 12 to 1
$wsIn.Range("B11").Value = 1

# maximumallowedaoutstandingbalance changes from 5000 to 10000
$wsIn.Range("B28").Value = 10000

# --- Append new account-mapping rows 31-42, copying formatting from row 13 ---
$wsIn.Range("A13:B13").Copy()
$wsIn.Range("A31:B42").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$newRows = @(
    @("fundsource", "Cash"),
    @("loanprotfolio", "Loan portfolio "),
    @("interestreceivable", "Interest Receivable "),
    @("penaltiesreceivable", "Penalties Receivable "),
    @("transferinsuspense", "Transfer in Suspence "),
    @("feesreceivable", "Fees Receivable"),
    @("incomefrominterest", "Income from interest"),
    @("incomefrompenalties", "Income from penalties"),
    @("incomefromfees", "Income from fees"),
    @("incomefromrecoveryrepayments", "Income from recovery repayments"),
    @("loseswrittenoff", "Losses Writtenoff "),
    @("overpaymentliability", "Overpayment Liability")
)

$r = 31
foreach ($pair in $newRows) {
    $wsIn.Cells.Item($r, 1).Value = $pair[0]
    $wsIn.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# --- ProductLoan_Output sheet edits ---
$wsOut.Range("B1").Value = "486-RBI-EPP-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"
$wsOut.Select() | Out-Null
$wsOut.Range("B1").Select() | Out-Null

# Re-select input sheet last so it stays the active/visible tab, scrolled
# down with the selection on B11 (matches the saved view state).
$wsIn.Select() | Out-Null
$wsIn.Range("B11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
